$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the values on row 5 to 2 decimal places (custom accuracy).
$ws.Range("B5").Value = 17.11
$ws.Range("C5").Value = 12.81
$ws.Range("D5").Value = 0.72
$ws.Range("E5").Value = 37.58
$ws.Range("F5").Value = 31.38
$ws.Range("G5").Value = 13.6
$ws.Range("H5").Value = 53.37
$ws.Range("I5").Value = 20.73
$ws.Range("J5").Value = 9.4
$ws.Range("K5").Value = 14.09
$ws.Range("L5").Value = 15.51
$ws.Range("M5").Value = 15.79
$ws.Range("N5").Value = 4.73
$ws.Range("O5").Value = 13.41
$ws.Range("P5").Value = 19.31
$ws.Range("Q5").Value = 11.13
$ws.Range("R5").Value = 0.16
$ws.Range("S5").Value = 0.68
$ws.Range("T5").Value = 198.48
$ws.Range("U5").Value = 37.52
$ws.Range("V5").Value = 12.39
$ws.Range("W5").Value = 25.48
$ws.Range("X5").Value = 13.28
$ws.Range("Y5").Value = 2.15
$ws.Range("Z5").Value = 25.86
$ws.Range("AA5").Value = 11.02
$ws.Range("AB5").Value = 9.55
$ws.Range("AC5").Value = 11.53
$ws.Range("AD5").Value = 16.27
$ws.Range("AE5").Value = 0.49
$ws.Range("AF5").Value = 48.89
$ws.Range("AG5").Value = 7.31
$ws.Range("AH5").Value = 15.43

# Remove the now-obsolete last data row (row 6), shifting the used range
# up so the sheet dimension becomes A1:AH5.
$ws.Rows.Item(6).Delete()
